$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(10, 2).Value = 47438
$ws.Cells.Item(10, 3).Value = "SIG-3w Lilliput LED Torch &amp; Table Lamp"
$ws.Cells.Item(10, 4).Value = 401.81
$ws.Cells.Item(10, 5).Value = 480.05
$ws.Cells.Item(10, 6).Value = 2
$ws.Cells.Item(10, 7).Value = 803.62
$ws.Cells.Item(11, 2).Value = 59408
$ws.Cells.Item(11, 3).Value = "SIG-3W Lilliput LED Torch &amp; Table Lamp"
$ws.Cells.Item(11, 4).Value = 388.17
$ws.Cells.Item(11, 5).Value = 463.78
$ws.Cells.Item(11, 6).Value = 9
$ws.Cells.Item(11, 7).Value = 3493.53
$ws.Cells.Item(26, 6).Value = 76
$ws.Cells.Item(26, 7).Value = 3503.6
$ws.Cells.Item(27, 6).Value = 49
$ws.Cells.Item(27, 7).Value = 1254.89
$ws.Cells.Item(32, 6).Value = 16
$ws.Cells.Item(32, 7).Value = 327.84
$ws.Cells.Item(46, 2).Value = 26241.49
$ws.Cells.Item(55, 6).Value = 207
$ws.Cells.Item(55, 7).Value = 39928.23
$ws.Cells.Item(65, 6).Value = 111
$ws.Cells.Item(65, 7).Value = 3409.92
$ws.Cells.Item(85, 2).Value = 153539.57
$ws.Cells.Item(95, 6).Value = 7
$ws.Cells.Item(95, 7).Value = 1758.75
$ws.Cells.Item(99, 6).Value = 16
$ws.Cells.Item(99, 7).Value = 1803.52
$ws.Cells.Item(103, 2).Value = 24048.32
$ws.Cells.Item(137, 6).Value = 50
$ws.Cells.Item(137, 7).Value = 4438
$ws.Cells.Item(147, 6).Value = 38
$ws.Cells.Item(147, 7).Value = 3846.74
$ws.Cells.Item(159, 2).Value = 69165.74000000001
$ws.Cells.Item(173, 6).Value = 3
$ws.Cells.Item(173, 7).Value = 126.96
$ws.Cells.Item(179, 6).Value = 68
$ws.Cells.Item(179, 7).Value = 3027.36
$ws.Cells.Item(180, 2).Value = 35484
$ws.Cells.Item(186, 6).Value = 7
$ws.Cells.Item(186, 7).Value = 117.74
$ws.Cells.Item(191, 6).Value = 56
$ws.Cells.Item(191, 7).Value = 7085.68
$ws.Cells.Item(193, 2).Value = 53925
$ws.Cells.Item(194, 2).Value = 57756
$ws.Cells.Item(198, 2).Value = 42957.2
$ws.Cells.Item(204, 6).Value = 10
$ws.Cells.Item(204, 7).Value = 7212.6
$ws.Cells.Item(205, 2).Value = 27824.7
$ws.Cells.Item(218, 6).Value = 23
$ws.Cells.Item(218, 7).Value = 1799.98
$ws.Cells.Item(228, 2).Value = 34629.71
$ws.Cells.Item(247, 6).Value = 2
$ws.Cells.Item(247, 7).Value = 62.86
$ws.Cells.Item(252, 6).Value = 41
$ws.Cells.Item(252, 7).Value = 2117.65
$ws.Cells.Item(259, 6).Value = 58
$ws.Cells.Item(259, 7).Value = 408.32
$ws.Cells.Item(267, 2).Value = 20410.49
$ws.Cells.Item(355, 6).Value = 32
$ws.Cells.Item(355, 7).Value = 5529.28
$ws.Cells.Item(361, 6).Value = 46
$ws.Cells.Item(361, 7).Value = 3390.66
$ws.Cells.Item(363, 6).Value = 35
$ws.Cells.Item(363, 7).Value = 2431.45
$ws.Cells.Item(365, 6).Value = 20
$ws.Cells.Item(365, 7).Value = 2636.6
$ws.Cells.Item(378, 6).Value = 22
$ws.Cells.Item(378, 7).Value = 2531.1
$ws.Cells.Item(379, 6).Value = 7
$ws.Cells.Item(379, 7).Value = 755.9299999999999
$ws.Cells.Item(381, 6).Value = 136
$ws.Cells.Item(381, 7).Value = 3144.32
$ws.Cells.Item(392, 6).Value = 2
$ws.Cells.Item(392, 7).Value = 171.44
$ws.Cells.Item(409, 6).Value = 150
$ws.Cells.Item(409, 7).Value = 25699.5
$ws.Cells.Item(418, 6).Value = 1
$ws.Cells.Item(418, 7).Value = 59.47
$ws.Cells.Item(419, 6).Value = 269
$ws.Cells.Item(419, 7).Value = 11077.42
$ws.Cells.Item(423, 2).Value = 136487.69
$ws.Cells.Item(486, 6).Value = 90
$ws.Cells.Item(486, 7).Value = 5463
$ws.Cells.Item(494, 6).Value = 84
$ws.Cells.Item(494, 7).Value = 4704.84
$ws.Cells.Item(497, 2).Value = 35552.28
$ws.Cells.Item(516, 6).Value = 178
$ws.Cells.Item(516, 7).Value = 12219.7
$ws.Cells.Item(520, 6).Value = 129
$ws.Cells.Item(520, 7).Value = 1268.07
$ws.Cells.Item(526, 6).Value = 747
$ws.Cells.Item(526, 7).Value = 72160.2
$ws.Cells.Item(527, 6).Value = 165
$ws.Cells.Item(527, 7).Value = 6139.65
$ws.Cells.Item(528, 6).Value = 189
$ws.Cells.Item(528, 7).Value = 4647.51
$ws.Cells.Item(529, 6).Value = 156
$ws.Cells.Item(529, 7).Value = 4194.84
$ws.Cells.Item(532, 2).Value = 150407.34
$ws.Cells.Item(563, 6).Value = 207
$ws.Cells.Item(563, 7).Value = 3336.84
$ws.Cells.Item(567, 2).Value = 50027.77
$ws.Cells.Item(592, 6).Value = 78
$ws.Cells.Item(592, 7).Value = 2441.4
$ws.Cells.Item(610, 2).Value = 54879.72
$ws.Cells.Item(618, 6).Value = 609
$ws.Cells.Item(618, 7).Value = 6528.48
$ws.Cells.Item(621, 6).Value = 228
$ws.Cells.Item(621, 7).Value = 13844.16
$ws.Cells.Item(623, 6).Value = 64
$ws.Cells.Item(623, 7).Value = 5488.64
$ws.Cells.Item(625, 6).Value = 19
$ws.Cells.Item(625, 7).Value = 1222.08
$ws.Cells.Item(638, 2).Value = 144778.51
$ws.Cells.Item(669, 6).Value = 46
$ws.Cells.Item(669, 7).Value = 3644.58
$ws.Cells.Item(672, 6).Value = 67
$ws.Cells.Item(672, 7).Value = 17828.03
$ws.Cells.Item(673, 6).Value = 140
$ws.Cells.Item(673, 7).Value = 3651.2
$ws.Cells.Item(678, 6).Value = 12
$ws.Cells.Item(678, 7).Value = 563.52
$ws.Cells.Item(688, 2).Value = 85652.16
$ws.Cells.Item(714, 6).Value = 60
$ws.Cells.Item(714, 7).Value = 10681.2
$ws.Cells.Item(715, 6).Value = 41
$ws.Cells.Item(715, 7).Value = 5352.55
$ws.Cells.Item(718, 6).Value = 108
$ws.Cells.Item(718, 7).Value = 2937.6
$ws.Cells.Item(719, 6).Value = 97
$ws.Cells.Item(719, 7).Value = 2638.4
$ws.Cells.Item(720, 2).Value = 26690.45
$ws.Cells.Item(745, 6).Value = 91
$ws.Cells.Item(745, 7).Value = 5632.9
$ws.Cells.Item(773, 2).Value = 133788.17
$ws.Cells.Item(780, 6).Value = 99
$ws.Cells.Item(780, 7).Value = 8464.5
$ws.Cells.Item(785, 2).Value = 14258.01
$ws.Cells.Item(821, 6).Value = 117
$ws.Cells.Item(821, 7).Value = 15572.7
$ws.Cells.Item(826, 6).Value = 325
$ws.Cells.Item(826, 7).Value = 22607
$ws.Cells.Item(829, 6).Value = 108
$ws.Cells.Item(829, 7).Value = 5802.84
$ws.Cells.Item(837, 2).Value = 190180.36
$ws.Cells.Item(840, 6).Value = 43
$ws.Cells.Item(840, 7).Value = 7429.11
$ws.Cells.Item(843, 6).Value = 63
$ws.Cells.Item(843, 7).Value = 6855.03
$ws.Cells.Item(856, 6).Value = 496
$ws.Cells.Item(856, 7).Value = 14939.52
$ws.Cells.Item(857, 6).Value = 364
$ws.Cells.Item(857, 7).Value = 28606.76
$ws.Cells.Item(858, 6).Value = 40
$ws.Cells.Item(858, 7).Value = 15112.4
$ws.Cells.Item(860, 6).Value = 38
$ws.Cells.Item(860, 7).Value = 19552.52
$ws.Cells.Item(863, 6).Value = 138
$ws.Cells.Item(863, 7).Value = 4389.78
$ws.Cells.Item(867, 2).Value = 197673.58
$ws.Cells.Item(891, 6).Value = 2
$ws.Cells.Item(891, 7).Value = 1147.12
$ws.Cells.Item(892, 6).Value = 1
$ws.Cells.Item(892, 7).Value = 573.5599999999999
$ws.Cells.Item(904, 2).Value = 36749.37
$ws.Cells.Item(923, 2).Value = 2524027.25
$ws.Cells.Item(924, 2).Value = 2524027.25